# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.609.54"
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("D3").Value = "1.958.24"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'244.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.84%  "
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").Value = "'58.70"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.34%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("D10").Value = "'56.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.81%  "
$ws.Range("D11").Value = "'0.0863"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +10.04%  "
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("D13").Value = "'22.33"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.48%  "
$ws.Range("E14").Value = "  -1.26%  "
$ws.Range("D15").Value = "2.244.60"
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("E16").Value = "  -1.34%  "
$ws.Range("E17").Value = "  -1.95%  "
$ws.Range("D18").Value = "1.959.33"
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("D19").Value = "36.508.52"
$ws.Range("E19").Value = "  +1.25%  "
$ws.Range("D20").Value = "0.0₃0880"
$ws.Range("E20").Value = "  +3.92%  "
$ws.Range("D21").Value = "'70.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.01%  "
$ws.Range("D22").Value = "'230.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E23").Value = "  -1.50%  "
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("E25").Value = "  -1.34%  "
$ws.Range("D26").Value = "'2.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.99%  "
$ws.Range("E27").Value = "  -1.61%  "
$ws.Range("D28").Value = "'162.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.09%  "
$ws.Range("D29").Value = "'0.137"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +10.15%  "
$ws.Range("D30").Value = "'19.68"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("E32").Value = "  +6.42%  "
$ws.Range("D33").Value = "'4.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.00%  "
$ws.Range("D34").Value = "'0.0646"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.02%  "
$ws.Range("E35").Value = "  -1.41%  "
$ws.Range("D36").Value = "'6.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.30%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "'2.20"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.70%  "
$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").Value = "'1.78"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.76%  "
$ws.Range("D40").Value = "'3.06"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.90%  "
$ws.Range("D41").Value = "'0.100"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.16%  "
$ws.Range("D42").Value = "'2.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.48%  "
$ws.Range("E43").Value = "  -2.15%  "
$ws.Range("E44").Value = "  +0.44%  "
$ws.Range("E45").Value = "  +2.88%  "
$ws.Range("E46").Value = "  -2.84%  "
$ws.Range("D47").Value = "1.357.97"
$ws.Range("E47").Value = "  +1.68%  "
$ws.Range("D48").Value = "'88.73"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.47%  "
$ws.Range("D49").Value = "'7.25"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.15%  "
$ws.Range("D50").Value = "'2.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("D51").Value = "'46.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.37%  "
